$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(0.06328177979961902, 0.3375848360084654, 0.1529057820181812, 6.48142807727062, 0, 7.035200475096886)
    3 = @(0.06328177979961902, 0.004309184025731883, 0.1529057820181812, 0.4998867070740569, 0, 0.720383452917589)
    4 = @(1.505614041169197, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 0, 4.371470058157054)
    5 = @(0.1554434735375247, 0.05231270169004087, 0.1529057820181812, 0.4998867070740569, 0, 0.8605486643198037)
    6 = @(3.182878228561681, 1.65323645889881, 3.082599426703578, 0.4998867070740569, 0, 8.418600821238126)
    7 = @(1.505614041169197, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 1, 4.371470058157054)
    8 = @(0.1554434735375247, 0.3375848360084654, 0.7127328510149897, 0.4998867070740569, 0, 1.705647867635037)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 6).Value = $vals[4]
    $ws.Cells.Item($row, 7).Value = $vals[5]
}
